$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 8362.375
$ws.Range("I69").Value = 1633.3334
$ws.Range("J69").Value = 12399.8
$ws.Range("K69").Value = 4900.0002
$ws.Range("L69").Value = 37199.39999999999
$ws.Range("M69").Value = -4026.0002
$ws.Range("N69").Value = -38947.39999999999
$ws.Range("H72").Value = 8362.375
$ws.Range("I72").Value = 1633.3334
$ws.Range("J72").Value = 12399.8
$ws.Range("K72").Value = 14700.0006
$ws.Range("L72").Value = 111598.2
$ws.Range("M72").Value = -10332.0006
$ws.Range("N72").Value = -120334.2
$ws.Range("H74").Value = 4312
$ws.Range("I74").Value = 4998
$ws.Range("J74").Value = 4214
$ws.Range("K74").Value = 4998
$ws.Range("L74").Value = 4214
$ws.Range("M74").Value = -4062
$ws.Range("N74").Value = -6086
$ws.Range("H76").Value = 3041.3794
$ws.Range("I76").Value = 3030.7693
$ws.Range("J76").Value = 3133.3333
$ws.Range("K76").Value = 3030.7693
$ws.Range("L76").Value = 3133.3333
$ws.Range("M76").Value = -2715.7693
$ws.Range("N76").Value = -3763.3333
$ws.Range("H77").Value = 4312
$ws.Range("I77").Value = 4998
$ws.Range("J77").Value = 4214
$ws.Range("K77").Value = 24990
$ws.Range("L77").Value = 21070
$ws.Range("M77").Value = -20310
$ws.Range("N77").Value = -30430
$ws.Range("H79").Value = 3041.3794
$ws.Range("I79").Value = 3030.7693
$ws.Range("J79").Value = 3133.3333
$ws.Range("K79").Value = 3030.7693
$ws.Range("L79").Value = 3133.3333
$ws.Range("M79").Value = -1938.7693
$ws.Range("N79").Value = -5317.3333
$ws.Range("H80").Value = 4800.8164
$ws.Range("I80").Value = 3785.2
$ws.Range("J80").Value = 5501.241
$ws.Range("K80").Value = 11355.6
$ws.Range("L80").Value = 16503.723
$ws.Range("M80").Value = -10357.6
$ws.Range("N80").Value = -18499.723
$ws.Range("H83").Value = 4800.8164
$ws.Range("I83").Value = 3785.2
$ws.Range("J83").Value = 5501.241
$ws.Range("K83").Value = 34066.8
$ws.Range("L83").Value = 49511.169
$ws.Range("M83").Value = -29074.8
$ws.Range("N83").Value = -59495.169
$ws.Range("H137").Value = 2611.625
$ws.Range("I137").Value = 1803.7273
$ws.Range("J137").Value = 4389
$ws.Range("K137").Value = 5411.1819
$ws.Range("L137").Value = 13167
$ws.Range("M137").Value = -2861.1819
$ws.Range("N137").Value = -18267

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3450.4106
$ws.Range("I32").Value = 2961.4897
$ws.Range("J32").Value = 6872.857
$ws.Range("K32").Value = 2961.4897
$ws.Range("L32").Value = 6872.857
$ws.Range("M32").Value = -2674.4897
$ws.Range("N32").Value = -7446.857
$ws.Range("H88").Value = 2348.4285
$ws.Range("I88").Value = 2177.25
$ws.Range("J88").Value = 2453.7693
$ws.Range("K88").Value = 2177.25
$ws.Range("L88").Value = 2453.7693
$ws.Range("M88").Value = -1771.25
$ws.Range("N88").Value = -3265.7693
$ws.Range("H91").Value = 2348.4285
$ws.Range("I91").Value = 2177.25
$ws.Range("J91").Value = 2453.7693
$ws.Range("K91").Value = 2177.25
$ws.Range("L91").Value = 2453.7693
$ws.Range("M91").Value = -773.25
$ws.Range("N91").Value = -5261.7693

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H55").Value = 29259.666
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 29259.666
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 29259.666
$ws.Range("N55").Value = -29805.666
$ws.Range("H81").Value = 22702.375
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 22702.375
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 22702.375
$ws.Range("N81").Value = -24824.375
$ws.Range("H84").Value = 22702.375
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 22702.375
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 68107.125
$ws.Range("N84").Value = -78715.125
$ws.Range("H86").Value = 1642.3684
$ws.Range("I86").Value = 1547
$ws.Range("J86").Value = 2000
$ws.Range("K86").Value = 1547
$ws.Range("L86").Value = 2000
$ws.Range("M86").Value = -424
$ws.Range("N86").Value = -4246
$ws.Range("H89").Value = 1642.3684
$ws.Range("I89").Value = 1547
$ws.Range("J89").Value = 2000
$ws.Range("K89").Value = 7735
$ws.Range("L89").Value = 10000
$ws.Range("M89").Value = -2119
$ws.Range("N89").Value = -21232
$ws.Range("H105").Value = 2175.4546
$ws.Range("I105").Value = 1607.5
$ws.Range("J105").Value = 2500
$ws.Range("K105").Value = 1607.5
$ws.Range("L105").Value = 2500
$ws.Range("M105").Value = 139.5
$ws.Range("N105").Value = -5994

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2755.3333
$ws.Range("I62").Value = 2268.182
$ws.Range("J62").Value = 3520.8572
$ws.Range("K62").Value = 2268.182
$ws.Range("L62").Value = 3520.8572
$ws.Range("M62").Value = -1644.182
$ws.Range("N62").Value = -4768.8572
$ws.Range("H65").Value = 2755.3333
$ws.Range("I65").Value = 2268.182
$ws.Range("J65").Value = 3520.8572
$ws.Range("K65").Value = 11340.91
$ws.Range("L65").Value = 17604.286
$ws.Range("M65").Value = -8220.91
$ws.Range("N65").Value = -23844.286

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 513.48334
$ws.Range("I113").Value = 526.64703
$ws.Range("J113").Value = 496.26923
$ws.Range("K113").Value = 1579.94109
$ws.Range("L113").Value = 1488.80769
$ws.Range("M113").Value = 590.0589100000002
$ws.Range("N113").Value = -5828.80769
$ws.Range("H131").Value = 888.8
$ws.Range("I131").Value = 818.25
$ws.Range("J131").Value = 894.93475
$ws.Range("K131").Value = 2454.75
$ws.Range("L131").Value = 2684.80425
$ws.Range("M131").Value = 2585.25
$ws.Range("N131").Value = -12764.80425

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H51").Value = 23360
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 23360
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 23360
$ws.Range("N51").Value = -24378
$ws.Range("H70").Value = 4883.7905
$ws.Range("I70").Value = 4799.2104
$ws.Range("J70").Value = 5017.7085
$ws.Range("K70").Value = 4799.2104
$ws.Range("L70").Value = 5017.7085
$ws.Range("M70").Value = -4529.2104
$ws.Range("N70").Value = -5557.7085
$ws.Range("H73").Value = 4883.7905
$ws.Range("I73").Value = 4799.2104
$ws.Range("J73").Value = 5017.7085
$ws.Range("K73").Value = 4799.2104
$ws.Range("L73").Value = 5017.7085
$ws.Range("M73").Value = -3863.2104
$ws.Range("N73").Value = -6889.7085
$ws.Range("H80").Value = 2843.8
$ws.Range("I80").Value = 2553.8462
$ws.Range("J80").Value = 3382.2856
$ws.Range("K80").Value = 2553.8462
$ws.Range("L80").Value = 3382.2856
$ws.Range("M80").Value = -1555.8462
$ws.Range("N80").Value = -5378.2856
$ws.Range("H83").Value = 2843.8
$ws.Range("I83").Value = 2553.8462
$ws.Range("J83").Value = 3382.2856
$ws.Range("K83").Value = 12769.231
$ws.Range("L83").Value = 16911.428
$ws.Range("M83").Value = -7777.231
$ws.Range("N83").Value = -26895.428

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 8641.177
$ws.Range("I68").Value = 12300
$ws.Range("J68").Value = 3414.2856
$ws.Range("K68").Value = 12300
$ws.Range("L68").Value = 3414.2856
$ws.Range("M68").Value = -11551
$ws.Range("N68").Value = -4912.2856
$ws.Range("H71").Value = 8641.177
$ws.Range("I71").Value = 12300
$ws.Range("J71").Value = 3414.2856
$ws.Range("K71").Value = 61500
$ws.Range("L71").Value = 17071.428
$ws.Range("M71").Value = -57756
$ws.Range("N71").Value = -24559.428

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 2611
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 2611
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 2611
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -3859
$ws.Range("H65").Value = 2611
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 2611
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 13055
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -19295
